$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '49.827.98'
$ws.Range('D2').Style = 'Normal'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.623.22'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +4.47%  '
$ws.Range('E3').Style = 'Normal'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('E4').Style = 'Normal'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '326.72'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.57%  '
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '109.94'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +1.26%  '
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.534'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +1.04%  '
$ws.Range('E7').Style = 'Normal'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E8').Style = 'Normal'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +3.37%  '
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '40.33'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +0.94%  '
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.67'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.58%  '
$ws.Range('E11').Style = 'Normal'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.03%  '
$ws.Range('E12').Style = 'Normal'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.27'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.22%  '
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.031.76'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +4.51%  '
$ws.Range('E15').Style = 'Normal'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.616.72'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +4.41%  '
$ws.Range('E16').Style = 'Normal'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +3.25%  '
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '49.754.57'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +3.37%  '
$ws.Range('E18').Style = 'Normal'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +11.33%  '
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.36'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.66%  '
$ws.Range('E20').Style = 'Normal'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.27%  '
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0₃0955'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +1.17%  '
$ws.Range('E22').Style = 'Normal'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '72.61'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +0.41%  '
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '278.05'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -0.29%  '
$ws.Range('E24').Style = 'Normal'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +1.09%  '
$ws.Range('E25').Style = 'Normal'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.52'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +2.80%  '
$ws.Range('E26').Style = 'Normal'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.00%  '
$ws.Range('E27').Style = 'Normal'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.23%  '
$ws.Range('E28').Style = 'Normal'
$ws.Range('B29').NumberFormat = '@'
$ws.Range('B29').Value = 'Cosmos'
$ws.Range('B29').Style = 'Normal'
$ws.Range('C29').NumberFormat = '@'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('C29').Style = 'Normal'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.96'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +1.65%  '
$ws.Range('E29').Style = 'Normal'
$ws.Range('B30').NumberFormat = '@'
$ws.Range('B30').Value = 'InjectiveProtocol'
$ws.Range('B30').Style = 'Normal'
$ws.Range('C30').NumberFormat = '@'
$ws.Range('C30').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('C30').Style = 'Normal'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '36.61'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +3.56%  '
$ws.Range('E30').Style = 'Normal'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +2.45%  '
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '49.75'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +0.40%  '
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.74'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -0.04%  '
$ws.Range('E33').Style = 'Normal'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +1.70%  '
$ws.Range('E34').Style = 'Normal'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.05%  '
$ws.Range('E35').Style = 'Normal'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.75%  '
$ws.Range('E36').Style = 'Normal'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.06'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +5.10%  '
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.74'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +1.49%  '
$ws.Range('E38').Style = 'Normal'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +5.87%  '
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '123.90'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +1.27%  '
$ws.Range('E40').Style = 'Normal'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.61%  '
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '22.50'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +4.79%  '
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.22'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.24%  '
$ws.Range('E43').Style = 'Normal'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +4.00%  '
$ws.Range('E44').Style = 'Normal'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +4.69%  '
$ws.Range('E45').Style = 'Normal'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.050.70'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +1.45%  '
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.30'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +15.20%  '
$ws.Range('E47').Style = 'Normal'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +8.51%  '
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.02'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.36%  '
$ws.Range('E49').Style = 'Normal'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +3.25%  '
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '81.52'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.94%  '
$ws.Range('E51').Style = 'Normal'
